$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new PR log row (row 18) matching the columns:
# ID | Title | Author | Approvers | Source / Target | Date
$row = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

$ws.Cells.Item($row, 1).Value = 42
$ws.Cells.Item($row, 2).Value = "Update README.md"
$ws.Cells.Item($row, 3).Value = "riya-morankar"
$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "edit1 to main"

# Force the Date cell to be stored as plain text (matching the rest of the
# column) instead of letting Excel auto-convert the "YYYY-MM-DD" string into
# a date serial number with a date number format.
$dateCell = $ws.Cells.Item($row, 6)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-06-18"
$dateCell.Style = "Normal"
